$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds the same date serial (45171) for every data
# row (2..289). The edit bumps that date by one day (45171 -> 45172) for
# every one of those rows, leaving everything else untouched.
$ws.Range("C2:C289").Value = 45172
